$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy header style from an existing header cell (e.g. AC1) to the new headers
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the Wins/Losses/Ties values for rows 2 through 64
for ($r = 2; $r -le 64; $r++) {
    $ws.Cells.Item($r, 30).Value = 91   # AD
    $ws.Cells.Item($r, 31).Value = 71   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
